$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated per-cluster ligand-side stats, keyed by "Sending cluster" (column A)
# Tuple order: LigandExpressingCells(E), LigandAvgExpr(G), LigandTotalExpr(H),
#              LigandSpecAvg(I), LigandSpecTotal(J)
$sendStats = @{
    "ECs"    = @(2, 1.0664485, 2.132897, 0.05919991215896408, 0.04426563970404137)
    "FAPs"   = @(3, 2.4368, 7.3104, 0.1352698662419833, 0.1517183119918233)
    "M1"     = @(3, 4.084561333333333, 12.253684, 0.2267391928829517, 0.2543100585687806)
    "M2"     = @(3, 4.804119, 14.412357, 0.2666827538331297, 0.2991106472783349)
    "Neutro" = @(3, 0.8298326666666666, 2.489498, 0.04606506640808776, 0.05166645248782833)
    "sCs"    = @(2, 4.792598, 9.585196, 0.2660432084748837, 0.1989288899691914)
}

# Updated per-cluster receptor-side stats, keyed by "Target cluster" (column D)
# Tuple order: ReceptorExpressingCells(K), ReceptorAvgExpr(M), ReceptorTotalExpr(N),
#              ReceptorSpecAvg(O), ReceptorSpecTotal(P)
$targetStats = @{
    "ECs"    = @(2, 7.5957635, 15.191527, 0.1268077702461478, 0.09750189942720215)
    "FAPs"   = @(3, 29.34977266666667, 88.049318, 0.4899809254318325, 0.5651160510901728)
    "M1"     = @(3, 0.3022816666666667, 0.906845, 0.005046453083523374, 0.00582029113900539)
    "M2"     = @(3, 0.7729826666666667, 2.318948, 0.01290458930151278, 0.01488341722809771)
    "Neutro" = @(3, 5.582807, 16.748421, 0.09320238938252692, 0.1074943196892874)
    "sCs"    = @(2, 16.2962195, 32.592439, 0.2720578725544566, 0.2091840214262345)
}

$lastRow = 37
for ($r = 2; $r -le $lastRow; $r++) {
    $sendKey = $ws.Cells.Item($r, 1).Value2
    $targetKey = $ws.Cells.Item($r, 4).Value2

    $s = $sendStats[$sendKey]
    $t = $targetStats[$targetKey]

    $E = $s[0]; $G = $s[1]; $H = $s[2]; $I = $s[3]; $J = $s[4]
    $K = $t[0]; $M = $t[1]; $N = $t[2]; $O = $t[3]; $P = $t[4]

    $Q = $G * $M
    $R = $H * $N
    $S = $I * $O
    $T = $J * $P

    $ws.Cells.Item($r, 5).Value  = $E   # E - Ligand-expressing cells
    $ws.Cells.Item($r, 7).Value  = $G   # G - Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $H   # H - Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $I   # I - Ligand derived specificity (avg)
    $ws.Cells.Item($r, 10).Value = $J   # J - Ligand derived specificity (total)
    $ws.Cells.Item($r, 11).Value = $K   # K - Receptor-expressing cells
    $ws.Cells.Item($r, 13).Value = $M   # M - Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $N   # N - Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $O   # O - Receptor derived specificity (avg)
    $ws.Cells.Item($r, 16).Value = $P   # P - Receptor derived specificity (total)
    $ws.Cells.Item($r, 17).Value = $Q   # Q - Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $R   # R - Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $S   # S - Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $T   # T - Edge total expression derived specificity
}
